$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = -12.90434384910118
$ws.Range("C4").Value = 4.440675708811348
$ws.Range("E6").Value = -4.72694934804071
$ws.Range("E7").Value = -4.442195584720943
$ws.Range("C8").Value = 1.600060471414855
$ws.Range("E8").Value = -2.94981949413492
$ws.Range("E9").Value = 2.83373437266663
$ws.Range("E10").Value = 0.9898525518332146
$ws.Range("C11").Value = 2.613530175870649
$ws.Range("E12").Value = 6.233648892986987
$ws.Range("C13").Value = -0.9738659311994247
$ws.Range("C14").Value = 1.725130460355073
$ws.Range("E14").Value = -1.194610791899997
$ws.Range("C15").Value = -2.221695202430862
$ws.Range("E15").Value = -6.324519326136457
$ws.Range("C16").Value = -1.290816600413769
$ws.Range("E16").Value = -8.277043919141525
$ws.Range("C17").Value = 2.066568225344989
$ws.Range("C18").Value = -1.187784794999081
$ws.Range("E19").Value = 8.305446157974018
